$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '52.106.36'
$ws.Range('E2').Value = '  +5.31%  '
$ws.Range('D3').Value = '2.792.06'
$ws.Range('E3').Value = '  +6.06%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '116.74'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +4.84%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '342.89'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +5.37%  '
$ws.Range('E7').Value = '  +4.09%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.579'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +6.52%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '42.10'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +6.79%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0869'
$ws.Range('D11').Style = 'Normal'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '20.13'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.16%  '
$ws.Range('E13').Value = '  +2.59%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.64'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +3.67%  '
$ws.Range('D15').Value = '3.234.67'
$ws.Range('E15').Value = '  +6.26%  '
$ws.Range('D16').Value = '2.793.20'
$ws.Range('E16').Value = '  +6.41%  '
$ws.Range('E17').Value = '  +4.31%  '
$ws.Range('D18').Value = '52.009.81'
$ws.Range('E18').Value = '  +5.18%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.23'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +11.40%  '
$ws.Range('E20').Value = '  +2.42%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.97'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +4.48%  '
$ws.Range('D22').Value = '0.0₃0991'
$ws.Range('E22').Value = '  +4.68%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '277.93'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +3.78%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '70.24'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.94%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.85'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +12.22%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '26.88'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +3.44%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.999'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.14%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.19'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.44%  '
$ws.Range('E29').Value = '  +1.14%  '
$ws.Range('E30').Value = '  +3.04%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '34.76'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.64%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '50.52'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.97%  '
$ws.Range('E33').Value = '  +4.01%  '
$ws.Range('E34').Value = '  +1.75%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.13'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +5.47%  '
$ws.Range('E36').Value = '  -0.08%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '18.96'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.25%  '
$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.98'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.67%  '
$ws.Range('B39').Value = 'LidoDAOToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.28'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +6.69%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.78'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +28.42%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0369'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +13.23%  '
$ws.Range('B42').Value = 'EnergySwap'
$ws.Range('C42').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '23.54'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +2.82%  '
$ws.Range('B43').Value = 'Monero'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '127.94'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.77%  '
$ws.Range('B44').Value = 'WEMIXToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.34'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +5.95%  '
$ws.Range('E45').Value = '  +3.62%  '
$ws.Range('D46').Value = '2.107.11'
$ws.Range('E46').Value = '  +3.21%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.33'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +3.28%  '
$ws.Range('E48').Value = '  +3.35%  '
$ws.Range('E49').Value = '  +6.73%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.919'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +22.28%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '8.94'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.19%  '
